$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every record row
# (rows 2 through 141). Update every one of these cells from 2023-09-03
# (serial 45172) to 2023-09-06 (serial 45175), keeping the rest of the
# workbook untouched.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
